$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring column H into existence as a sibling of column G, inheriting G's
# formatting (bold header font, border, centered alignment) the same way
# Excel extends formatting when a new column is inserted right after an
# existing one.
$ws.Range("H1").EntireColumn.Insert()

# New "Save" column: header + per-row flag values.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
